$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price(D) and Volume(E) columns to text format so numeric-looking
# strings (e.g. "1.014", "28.107.65") are stored as text, matching the
# original inlineStr cell type instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.107.65"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").Value = "1.891.32"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +1.08%  "

$ws.Range("D5").Value = "337.05"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("E6").Value = "  +1.01%  "

$ws.Range("D7").Value = "0.4758"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").Value = "0.3963"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").Value = "0.08056"
$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").Value = "1.023"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "22.02"
$ws.Range("E12").Value = "  +1.10%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "6.047"
$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.872.78"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "7.245"
$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("D17").Value = "88.70"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("D18").Value = "0.06774"
$ws.Range("E18").Value = "  +2.06%  "

$ws.Range("D19").Value = "0.00001055"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "17.13"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").Value = "28.084.75"
$ws.Range("E22").Value = "  +1.42%  "

$ws.Range("D23").Value = "5.543"
$ws.Range("E23").Value = "  +0.97%  "

$ws.Range("D24").Value = "11.06"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D26").Value = "2.097.63"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").Value = "160.96"
$ws.Range("E27").Value = "  +1.59%  "

$ws.Range("D28").Value = "20.07"
$ws.Range("E28").Value = "  -0.80%  "

$ws.Range("D29").Value = "2.116"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").Value = "5.557"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").Value = "122.24"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").Value = "0.9819"
$ws.Range("E32").Value = "  +1.74%  "

$ws.Range("D33").Value = "0.09614"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("D34").Value = "3.641"
$ws.Range("E34").Value = "  +1.36%  "

$ws.Range("D35").Value = "5.376"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("D36").Value = "1.372"
$ws.Range("E36").Value = "  -5.01%  "

$ws.Range("D37").Value = "0.02260"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "0.06097"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("D39").Value = "1.206"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").Value = "8.226"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").Value = "1.012"
$ws.Range("E41").Value = "  +1.10%  "

$ws.Range("D42").Value = "0.5990"
$ws.Range("E42").Value = "  +0.42%  "

$ws.Range("D43").Value = "0.1900"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").Value = "10.39"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("D45").Value = "1.267"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").Value = "0.5684"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "12.18"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("D48").Value = "1.939"
$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("D49").Value = "3.368"
$ws.Range("E49").Value = "  -0.57%  "

$ws.Range("D50").Value = "0.06839"
$ws.Range("E50").Value = "  -0.22%  "

$ws.Range("D51").Value = "112.58"
$ws.Range("E51").Value = "  -1.46%  "

# Restore the default style so the cells keep the original (no explicit
# style index) appearance rather than the Text-formatted style.
$ws.Range("D2:E51").Style = "Normal"
